# Apply "Add data for 2021-12-26" update to the carjacking-by-neighborhood
# workbook: rename the sheet / header to reflect data through 2021-12-18,
# and update the affected neighborhood/month cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet name + header label: December 17 -> December 18 ---
$ws.Name = "Through 2021-12-18"
$ws.Range("B1").Value = "December 2021 (through December 18)"

# --- Updated counts (existing cells) ---
$ws.Range("Z4").Value = 4    # North Lawndale   / December 2019
$ws.Range("N7").Value = 6    # Austin           / December 2020
$ws.Range("Z7").Value = 3    # Austin           / December 2019
$ws.Range("AL7").Value = 9   # Austin           / December 2018
$ws.Range("BJ8").Value = 4   # Chatham          / December 2016
$ws.Range("B9").Value = 6    # Grand Crossing   / December 2021 (through Dec 18)
$ws.Range("AX18").Value = 3  # Grand Boulevard  / December 2017
$ws.Range("B21").Value = 4   # Wicker Park      / December 2021 (through Dec 18)
$ws.Range("N40").Value = 4   # Calumet Heights  / December 2020
$ws.Range("AL66").Value = 3  # Avondale         / December 2018
$ws.Range("B70").Value = 3   # Bucktown         / December 2021 (through Dec 18)

# --- New counts (previously-empty cells) ---
$ws.Range("AX9").Value = 1   # Grand Crossing   / December 2017
$ws.Range("BJ10").Value = 1  # Douglas          / December 2016
$ws.Range("AL13").Value = 1  # Roseland         / December 2018
$ws.Range("BJ21").Value = 1  # Wicker Park      / December 2016
$ws.Range("BV39").Value = 1  # Brighton Park    / December 2015
$ws.Range("AX41").Value = 1  # Chinatown        / December 2017
$ws.Range("BV62").Value = 1  # Armour Square    / December 2015
$ws.Range("N99").Value = 1   # West Lawn        / December 2020
